$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the SuiNetwork/Fetch.AI row-order swap)

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '59.007.67'
$ws.Cells.Item(2, 5).Value = '  +0.86%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.572.76'
$ws.Cells.Item(3, 5).Value = '  -0.23%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '562.84'
$ws.Cells.Item(5, 5).Value = '  +3.85%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.31'
$ws.Cells.Item(6, 5).Value = '  -1.05%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.06%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.593'
$ws.Cells.Item(8, 5).Value = '  +1.76%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.577.32'
$ws.Cells.Item(9, 5).Value = '  -0.27%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.63'
$ws.Cells.Item(10, 5).Value = '  -1.76%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +2.31%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.151'
$ws.Cells.Item(12, 5).Value = '  +8.79%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.35%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '59.091.19'
$ws.Cells.Item(15, 5).Value = '  +1.11%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '21.80'
$ws.Cells.Item(16, 5).Value = '  +6.09%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000136'
$ws.Cells.Item(17, 5).Value = '  +3.30%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.582.54'
$ws.Cells.Item(18, 5).Value = '  +0.79%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.48'
$ws.Cells.Item(19, 5).Value = '  +0.85%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '334.88'
$ws.Cells.Item(20, 5).Value = '  +0.37%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.14'
$ws.Cells.Item(21, 5).Value = '  +1.03%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +1.38%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.04%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.57'
$ws.Cells.Item(24, 5).Value = '  -2.75%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +5.30%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.22%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.161'
$ws.Cells.Item(27, 5).Value = '  +2.03%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.19'
$ws.Cells.Item(28, 5).Value = '  +1.82%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.0₃0779'
$ws.Cells.Item(29, 5).Value = '  +5.53%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.998'
$ws.Cells.Item(30, 5).Value = '  -0.03%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +2.58%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '159.84'
$ws.Cells.Item(32, 5).Value = '  +4.52%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.03'
$ws.Cells.Item(33, 5).Value = '  +0.85%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '18.84'
$ws.Cells.Item(34, 5).Value = '  -0.33%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.42%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.881'
$ws.Cells.Item(36, 5).Value = '  +7.49%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'SuiNetwork'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.875'
$ws.Cells.Item(37, 5).Value = '  +3.51%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +2.92%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '36.67'
$ws.Cells.Item(39, 5).Value = '  -1.10%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +3.99%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '295.78'
$ws.Cells.Item(41, 5).Value = '  +6.25%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.61'

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.03%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0970'
$ws.Cells.Item(44, 5).Value = '  +2.96%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.592'
$ws.Cells.Item(45, 5).Value = '  +0.22%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.97%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.60'
$ws.Cells.Item(47, 5).Value = '  -0.26%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '125.35'
$ws.Cells.Item(48, 5).Value = '  +14.89%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '18.93'
$ws.Cells.Item(49, 5).Value = '  +2.38%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0231'
$ws.Cells.Item(50, 5).Value = '  +1.91%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '18.34'
$ws.Cells.Item(51, 5).Value = '  +2.79%  '
